$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Status paragraph (para 3): describe the 2023 update, and the prior
#    "beta" -> "final" (Aug 21, 2022) milestone, instead of just saying
#    "final release ... first released in beta form on September 9, 2021."
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(3).Range
$p1.Find.Execute(
    "final release of Indigo Book 2.0, which was first released in beta form on September 9, 2021. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2023 update of Indigo Book 2.0, which was first released in beta form on September 9, 2021; and then in final form on August 21, 2022. ",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Publisher paragraph (para 5): re-key the same sentence (no wording
#    change) so the run spanning "Public.Resource.Org" is unified with its
#    neighbors. (Stop just before the trailing non-breaking space so it is
#    not clobbered by a literal space.)
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(5).Range
$p2.Find.Execute(
    "This file was published by Public.Resource.Org, Inc., (“Public Resource”) a California nonprofit corporation registered under I.R.C. § 501(c)(3). Contact information for Public Resource is at",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This file was published by Public.Resource.Org, Inc., (“Public Resource”) a California nonprofit corporation registered under I.R.C. § 501(c)(3). Contact information for Public Resource is at",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3) CC-0 paragraph (para 6): re-key "—"No Rights Reserved"..." (no wording
#    change) to unify the run split.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(6).Range
$p3.Find.Execute(
    "—“No Rights Reserved” and we waive all copyright and related rights in this work.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "—“No Rights Reserved” and we waive all copyright and related rights in this work.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 4) "Statement of Nonaffiliation" heading (para 9): re-key to unify the
#    spell-check-split run.
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs(9).Range
$p4.Find.Execute(
    "Statement of Nonaffiliation",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Statement of Nonaffiliation",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 5) "Cite as" paragraph (para 12): append a sentence about the 2023 update
#    of the Second Edition, after the existing citation.
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs(12).Range
$p5.Find.Execute(
    "Public.Resource.Org 2d ed. 2021).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Public.Resource.Org 2d ed. 2021). Where relevant, the citation may also indicate that this is the 2023 update of the Second Edition.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 6) Formats paragraph (para 14): "HTML AND PDF" -> "HTML and PDF"
# ---------------------------------------------------------------------------
$p6 = $d.Paragraphs(14).Range
$p6.Find.Execute(
    "This document is available in HTML AND PDF formats.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This document is available in HTML and PDF formats.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 7) Bluepages paragraph (para 21): re-key "...'s "Bluepages"—that is," to
#    unify the spell-check-split run. (Stop just before the trailing
#    non-breaking space so it is not clobbered by a literal space.)
# ---------------------------------------------------------------------------
$p7 = $d.Paragraphs(21).Range
$p7.Find.Execute(
    "’s “Bluepages”—that is,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "’s “Bluepages”—that is,",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 8) NYU / Sprigman paragraph (para 25): re-key "...Christopher Jon
#    Sprigman" to unify the spell-check-split run.
# ---------------------------------------------------------------------------
$p8 = $d.Paragraphs(25).Range
$p8.Find.Execute(
    "working under the direction of Professor Christopher Jon Sprigman",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "working under the direction of Professor Christopher Jon Sprigman",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 9) Same paragraph (para 25): re-key "...compiled by Professor Jennifer
#    Murphy Romig with assistance..." to unify the spell-check-split run.
# ---------------------------------------------------------------------------
$p9 = $d.Paragraphs(25).Range
$p9.Find.Execute(
    "was compiled by Professor Jennifer Murphy Romig with assistance",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "was compiled by Professor Jennifer Murphy Romig with assistance",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 10) Signature block (para 33): re-key "Jennifer Murphy Romig" to unify the
#     spell-check-split run.
# ---------------------------------------------------------------------------
$p10 = $d.Paragraphs(33).Range
$p10.Find.Execute(
    "Professor Jennifer Murphy Romig",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Professor Jennifer Murphy Romig",
    2) | Out-Null
